$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "Columns:" helper list in column B (originally B60:B123) is the source
# list of selectable screener column names. This commit adds a batch of new
# technical-indicator names into that already-alphabetised list:
#   ADX, EMA10, EMA20, EMA30, EMA50, EMA100, EMA200, Mom, RSI,
#   SMA10, SMA20, SMA30, SMA50, SMA100, SMA200, Stoch.K, Stoch.D,
#   StochRSI.K, UO, VWMA
# which pushes the existing BBP..WillR tail (and everything after it) down
# by 20 rows. Insert 20 blank rows right before the old "BBP" row (83) and
# fill the new block (now rows 83-109) with the full merged, sorted list.
# ---------------------------------------------------------------------------

$ws.Range("A83:A102").EntireRow.Insert()

$newColumnNames = @(
    "ADX",
    "BBP",
    "CCI",
    "EMA10",
    "EMA20",
    "EMA30",
    "EMA50",
    "EMA100",
    "EMA200",
    "HMA",
    "Ichimoku",
    "MACD",
    "MACD.Signal",
    "Mom",
    "RSI",
    "SMA10",
    "SMA20",
    "SMA30",
    "SMA50",
    "SMA100",
    "SMA200",
    "Stoch.K",
    "Stoch.D",
    "StochRSI.K",
    "UO",
    "VWMA",
    "WillR"
)

$startRow = 83
for ($i = 0; $i -lt $newColumnNames.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $newColumnNames[$i]
}
